$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear existing data rows (rows 2-18) before writing new dataset (rows 2-44)
$ws.Range("A2:D18").Clear()

$data = @(
    ,@("rfomcalp", 7366, 958.8200000000002, 460.9415740180638)
    ,@("z", 30, 3.32, 542.1686746987951)
    ,@("ashmom", 1409, 155.58, 543.3860393366757)
    ,@("RFOCHAN", 274, 32.07, 512.6286248830683)
    ,@("rfoaves", 891, 99.8, 535.6713426853707)
    ,@("rfoschwa", 1975, 323.77, 366.0005559502115)
    ,@("rfolindsj", 2923, 441.88, 396.8950846383634)
    ,@("rfoalves", 2563, 281.02, 547.2208383744929)
    ,@("RFORoot", 91, 36.06, 151.414309484193)
    ,@("RFOCARSOS", 6537, 877.55, 446.948891801037)
    ,@("rfodel", 7, 0.22, 1909.090909090909)
    ,@("RFODISCH", 3743, 538.17, 417.3030826690451)
    ,@("rfosilvet", 3033, 281.51, 646.4423999147455)
    ,@("RFOKSB162", 2295, 186.09, 739.9645332903434)
    ,@("MAZZA", 2110, 190.27, 665.3702633100331)
    ,@("rfopiercn", 1132, 105.59, 643.2427313192537)
    ,@("rfoksc071", 3008, 325.04, 555.254737878415)
    ,@("rforootj", 2419, 224.01, 647.9175036828713)
    ,@("rfoghart", 1086, 68.90000000000001, 945.7184325108852)
    ,@("rfoksc905", 582, 59.72, 584.7287340924314)
    ,@("RFOROOT", 1711, 217.03, 473.0221628346311)
    ,@("rfoander", 1763, 201.78, 524.2343146000595)
    ,@("RFOKSB167", 1853, 186.01, 597.7098005483576)
    ,@("rfoely", 90, 7.42, 727.7628032345013)
    ,@("rfoalve", 1523, 298.83, 305.7925911053107)
    ,@("RFOBRADA", 2528, 190.04, 798.1477583666597)
    ,@("rfoksb170", 440, 58.35, 452.4421593830334)
    ,@("rfalves", 79, 26.81, 176.7997016038792)
    ,@("RFOKSB699", 1756, 181.98, 578.9647213979558)
    ,@("rfodeloy", 7377, 880.4299999999999, 502.7316197767)
    ,@("RFOKSB164", 1707, 188.7, 542.7662957074722)
    ,@("RFOCLARK", 2020, 226.74, 534.5329452236041)
    ,@("long", 1073, 133.31, 482.9345135398694)
    ,@("RFOKS171", 1900, 188.24, 605.6098597535062)
    ,@("mh", 3993, 483.08, 495.9427010019044)
    ,@("RFOCOOK", 1729, 191.31, 542.2612513721185)
    ,@("RFOBROCK", 7023, 779.0399999999998, 540.8964879852126)
    ,@(".", 1933, 188.63, 614.8544770185018)
    ,@("RFOKSB419", 911, 104.85, 521.3161659513592)
    ,@("rfoksc041", 7727, 951.77, 487.1134832995367)
    ,@("rfogatli", 1291, 113.89, 680.1299499517077)
    ,@("RFOKSB183", 182, 20.8, 525)
    ,@("lee", 81, 5.63, 863.2326820603907)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowNum = $i + 2
    $rowData = $data[$i]
    $ws.Cells.Item($rowNum, 1).Value = $rowData[0]
    $ws.Cells.Item($rowNum, 2).Value = $rowData[1]
    $ws.Cells.Item($rowNum, 3).Value = $rowData[2]
    $ws.Cells.Item($rowNum, 4).Value = $rowData[3]
}
